$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 / J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").Borders.LineStyle = 1
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160

# Data for I2:J69
$data = @(
    @(5,6),
    @(8,8),
    @(6,7),
    @(5,6),
    @(6,6),
    @(9,9),
    @(9,9),
    @(6,6),
    @(6,7),
    @(7,7),
    @(5,5),
    @(5,6),
    @(8,8),
    @(6,7),
    @(7,7),
    @(5,6),
    @(6,6),
    @(8,8),
    @(7,7),
    @(7,7),
    @(6,7),
    @(8,8),
    @(7,7),
    @(8,8),
    @(6,6),
    @(5,6),
    @(9,9),
    @(11,11),
    @(8,8),
    @(6,7),
    @(6,7),
    @(6,6),
    @(8,8),
    @(3,5),
    @(8,8),
    @(6,7),
    @(7,7),
    @(6,6),
    @(7,7),
    @(7,7),
    @(5,5),
    @(9,10),
    @(8,8),
    @(12,12),
    @(7,8),
    @(7,8),
    @(8,9),
    @(7,7),
    @(10,10),
    @(6,7),
    @(8,8),
    @(7,8),
    @(7,7),
    @(6,6),
    @(9,9),
    @(6,7),
    @(9,9),
    @(8,8),
    @(7,7),
    @(8,8),
    @(8,8),
    @(8,8),
    @(7,8),
    @(7,7),
    @(7,7),
    @(9,9),
    @(7,7),
    @(3,3)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $pair = $data[$idx]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
